$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.910.45'
$ws.Cells.Item(2, 5).Value = '  -1.07%  '

$ws.Cells.Item(3, 4).Value = '2.238.24'
$ws.Cells.Item(3, 5).Value = '  -2.01%  '

$ws.Cells.Item(4, 5).Value = '  +0.11%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '113.16'
$ws.Cells.Item(5, 4).NumberFormat = 'General'
$ws.Cells.Item(5, 5).Value = '  +0.31%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '278.63'
$ws.Cells.Item(6, 4).NumberFormat = 'General'
$ws.Cells.Item(6, 5).Value = '  +5.15%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.628'
$ws.Cells.Item(7, 4).NumberFormat = 'General'
$ws.Cells.Item(7, 5).Value = '  +0.18%  '

$ws.Cells.Item(8, 5).Value = '  +0.16%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.607'
$ws.Cells.Item(9, 4).NumberFormat = 'General'
$ws.Cells.Item(9, 5).Value = '  -0.27%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '46.11'
$ws.Cells.Item(10, 4).NumberFormat = 'General'
$ws.Cells.Item(10, 5).Value = '  -2.03%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0929'
$ws.Cells.Item(11, 4).NumberFormat = 'General'
$ws.Cells.Item(11, 5).Value = '  -0.79%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '9.02'
$ws.Cells.Item(12, 4).NumberFormat = 'General'
$ws.Cells.Item(12, 5).Value = '  -3.14%  '

$ws.Cells.Item(13, 5).Value = '  -3.01%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '15.29'
$ws.Cells.Item(14, 4).NumberFormat = 'General'
$ws.Cells.Item(14, 5).Value = '  -1.26%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.870'
$ws.Cells.Item(15, 4).NumberFormat = 'General'
$ws.Cells.Item(15, 5).Value = '  +0.47%  '

$ws.Cells.Item(16, 4).Value = '2.575.69'
$ws.Cells.Item(16, 5).Value = '  -2.08%  '

$ws.Cells.Item(17, 4).Value = '2.244.08'
$ws.Cells.Item(17, 5).Value = '  -1.93%  '

$ws.Cells.Item(18, 4).Value = '42.936.49'
$ws.Cells.Item(18, 5).Value = '  -0.80%  '

$ws.Cells.Item(19, 5).Value = '  -1.13%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.75'
$ws.Cells.Item(20, 4).NumberFormat = 'General'
$ws.Cells.Item(20, 5).Value = '  -0.63%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '72.03'
$ws.Cells.Item(21, 4).NumberFormat = 'General'
$ws.Cells.Item(21, 5).Value = '  +0.10%  '

$ws.Cells.Item(22, 5).Value = '  -5.22%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.03'
$ws.Cells.Item(23, 4).NumberFormat = 'General'
$ws.Cells.Item(23, 5).Value = '  +5.99%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '231.47'
$ws.Cells.Item(24, 4).NumberFormat = 'General'
$ws.Cells.Item(24, 5).Value = '  -1.53%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '9.25'
$ws.Cells.Item(25, 4).NumberFormat = 'General'
$ws.Cells.Item(25, 5).Value = '  -2.76%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '12.15'
$ws.Cells.Item(26, 4).NumberFormat = 'General'
$ws.Cells.Item(26, 5).Value = '  +6.49%  '

$ws.Cells.Item(27, 5).Value = '  -0.95%  '

$ws.Cells.Item(28, 5).Value = '  -1.78%  '

$ws.Cells.Item(29, 5).Value = '  -0.29%  '

$ws.Cells.Item(30, 5).Value = '  -2.66%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '173.53'
$ws.Cells.Item(31, 4).NumberFormat = 'General'
$ws.Cells.Item(31, 5).Value = '  +0.21%  '

$ws.Cells.Item(32, 5).Value = '  -2.24%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.0903'
$ws.Cells.Item(33, 4).NumberFormat = 'General'
$ws.Cells.Item(33, 5).Value = '  -0.35%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '5.57'
$ws.Cells.Item(34, 4).NumberFormat = 'General'
$ws.Cells.Item(34, 5).Value = '  -2.31%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '4.30'
$ws.Cells.Item(35, 4).NumberFormat = 'General'
$ws.Cells.Item(35, 5).Value = '  +5.93%  '

$ws.Cells.Item(36, 5).Value = '  -0.40%  '

$ws.Cells.Item(37, 5).Value = '  +0.19%  '

$ws.Cells.Item(38, 5).Value = '  +0.43%  '

$ws.Cells.Item(39, 5).Value = '  +1.54%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.57'
$ws.Cells.Item(40, 4).NumberFormat = 'General'
$ws.Cells.Item(40, 5).Value = '  -1.99%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '70.91'
$ws.Cells.Item(41, 4).NumberFormat = 'General'
$ws.Cells.Item(41, 5).Value = '  -7.04%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '13.18'
$ws.Cells.Item(42, 4).NumberFormat = 'General'
$ws.Cells.Item(42, 5).Value = '  -6.18%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.232'
$ws.Cells.Item(43, 4).NumberFormat = 'General'
$ws.Cells.Item(43, 5).Value = '  -2.68%  '

$ws.Cells.Item(44, 5).Value = '  -0.05%  '

$ws.Cells.Item(45, 5).Value = '  -3.56%  '

$ws.Cells.Item(46, 5).Value = '  -8.10%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.27'
$ws.Cells.Item(47, 4).NumberFormat = 'General'
$ws.Cells.Item(47, 5).Value = '  +0.68%  '

$ws.Cells.Item(48, 5).Value = '  -2.64%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.0986'
$ws.Cells.Item(49, 4).NumberFormat = 'General'
$ws.Cells.Item(49, 5).Value = '  -1.05%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '100.36'
$ws.Cells.Item(50, 4).NumberFormat = 'General'
$ws.Cells.Item(50, 5).Value = '  -3.78%  '

$ws.Cells.Item(51, 2).Value = 'TheSandbox'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.636'
$ws.Cells.Item(51, 4).NumberFormat = 'General'
$ws.Cells.Item(51, 5).Value = '  +5.72%  '
